$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Range("A6:B7").Insert(-4121)
